# Wijzigingen en aanvullingen tijdreeksen.
#
# - Rename the only worksheet from "Blad1" to "data"
# - Move the active cell selection from B19 to A18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "data"
$ws.Range("A18").Select()
